$d = $word.ActiveDocument

# The "ColorCard" bullet currently ends with a placeholder "(...)". Flesh it
# out with the actual enumeration of supported video card codes.
$colorCardRange = $d.Content
$replaced = $colorCardRange.Find.Execute( `
    "(...)", $true, $false, $false, $false, $false, $true, 1, $false, `
    "(3 = VGA, 4 = EGA, 6 = Tandy, 7 = CGA/PCjr, 10 = Hercules)", 2)

# Word tracks the location of the user's last edit with the hidden "_GoBack"
# bookmark. Since the edit above lands on the ColorCard line, move that
# bookmark there too -- right before the closing parenthesis we just typed.
$locateRange = $d.Content
$located = $locateRange.Find.Execute( `
    "Hercules)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($located) {
    $gobackPos = $locateRange.End - 1
    $gobackRange = $d.Range($gobackPos, $gobackPos)
    $d.Bookmarks.Add("_GoBack", $gobackRange)
}

Write-Host ("ColorCard text updated: " + $replaced + "; _GoBack relocated: " + $located)
